$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D" = 0.0509
    "E" = 0.141
    "G" = 0.6349809885931559
    "H" = 0.6349809885931559
    "I" = 0.4638783269961977
    "J" = 0.2935938778456948
    "K" = 2.5
    "L" = 0.4752851711026616
    "M" = 1.52
    "N" = 0.08397790055248618
    "O" = 0.608
    "P" = 1.52
    "Q" = 0.08397790055248618
    "R" = 0.608
    "U" = 0.305
    "V" = 0.01685082872928177
    "W" = 0.1724137931034483
    "X" = 0.07422547200288908
    "Y" = 0.09818832110055921
    "Z" = 0.3750445632798574
    "AA" = 0.1101107876982784
    "AB" = 0.07422547200288908
    "AC" = 0.03588531569538932
    "AG" = -0.305
    "AJ" = -0.01713964596796853
    "AK" = -0.02020536601523683
    "AM" = -1.47
    "AP" = -0.1013289036544851
    "AQ" = -1.659863945578231
}

foreach ($row in 2,3) {
    foreach ($col in $values.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $values[$col]
    }
}
